$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2; "B"="1.02"; "C"="1.04017524383956"; "D"="1.042819636062604"; "E"="1.047274954924097"; "F"="1.0555957671663"; "I"="1.02359499962809"; "J"="1.045263271193739"; "K"="1.045595246361183"; "L"="1.050038048749578"; "M"="1.058335791655322"; "N"="1.0183490914852" }
    @{ Row=3; "B"="1.02"; "C"="1.042716868103572"; "D"="1.045276176802779"; "E"="1.049600976006424"; "F"="1.058140979816293"; "I"="1.023504579208684"; "J"="1.047442465382062"; "K"="1.047858394805339"; "L"="1.052171958826295"; "M"="1.060690067838607"; "N"="1.019131940075213" }
    @{ Row=4; "B"="1.02"; "C"="1.044352909264869"; "D"="1.046857622668442"; "E"="1.051098044896564"; "F"="1.059779841418709"; "I"="1.023443335729026"; "J"="1.048844114537246"; "K"="1.049314436128667"; "L"="1.053544418375484"; "M"="1.062205120318036"; "N"="1.019634041792785" }
    @{ Row=5; "B"="1.02"; "C"="1.045038707707896"; "D"="1.04752057609493"; "E"="1.051725542078274"; "F"="1.060466943123821"; "I"="1.023416932628352"; "J"="1.049431397062961"; "K"="1.049924602763515"; "L"="1.054119455211515"; "M"="1.06284010745029"; "N"="1.019844076417721" }
    @{ Row=6; "B"="1.02"; "C"="1.045153741009738"; "D"="1.047631779792948"; "E"="1.051830793363836"; "F"="1.060582202090532"; "I"="1.023412460900047"; "J"="1.049529890216213"; "K"="1.050026939332133"; "L"="1.054215893725982"; "M"="1.062946612243976"; "N"="1.019879281142692" }
    @{ Row=7; "B"="1.02"; "C"="1.044362080694967"; "D"="1.046866488426635"; "E"="1.051106436828207"; "F"="1.059789029803537"; "I"="1.02344298551034"; "J"="1.048851969514203"; "K"="1.049322596813128"; "L"="1.053552109625919"; "M"="1.062213612605454"; "N"="1.01963685238427" }
    @{ Row=8; "B"="1.02"; "C"="1.041036003972119"; "D"="1.043651544782026"; "E"="1.048062736051536"; "F"="1.056457635697644"; "I"="1.023565007130014"; "J"="1.046001516165925"; "K"="1.046361850760844"; "L"="1.050760965744564"; "M"="1.059133184507377"; "N"="1.018614591491492" }
    @{ Row=9; "B"="1.02"; "C"="1.035106927191074"; "D"="1.037921960566549"; "E"="1.042635629575311"; "F"="1.050523089021407"; "I"="1.023759139049229"; "J"="1.04091183494964"; "K"="1.041078281382405"; "L"="1.04577671194103"; "M"="1.053639007243571"; "N"="1.01677833938926" }
    @{ Row=10; "B"="1.02"; "C"="1.031104808891264"; "D"="1.03405550218697"; "E"="1.038971493325861"; "F"="1.046520062579966"; "I"="1.023874603872166"; "J"="1.037470649832529"; "K"="1.037508062639323"; "L"="1.042406508592142"; "M"="1.049928490539339"; "N"="1.015529589008377" }
    @{ Row=11; "B"="1.02"; "C"="1.029359312437535"; "D"="1.032369425976712"; "E"="1.03737322502341"; "F"="1.044774848800066"; "I"="1.023921309042974"; "J"="1.03596846709841"; "K"="1.03595004106279"; "L"="1.040935245259231"; "M"="1.048309724469146"; "N"="1.014982781596006" }
    @{ Row=12; "B"="1.02"; "C"="1.028709000077692"; "D"="1.031741290895359"; "E"="1.036777739503772"; "F"="1.044124745288904"; "I"="1.023938164268118"; "J"="1.035408605547208"; "K"="1.035369441849503"; "L"="1.040386897755595"; "M"="1.047706561181219"; "N"="1.014778735749135" }
    @{ Row=13; "B"="1.02"; "C"="1.028848584021226"; "D"="1.031876112877348"; "E"="1.036905556456137"; "F"="1.044264279696542"; "I"="1.023934571058996"; "J"="1.03552878401884"; "K"="1.035494068494841"; "L"="1.040504605077418"; "M"="1.047836027963775"; "N"="1.014822547050643" }
    @{ Row=14; "B"="1.02"; "C"="1.029305597761362"; "D"="1.032317542296349"; "E"="1.037324039455552"; "F"="1.044721149255309"; "I"="1.023922712353277"; "J"="1.035922227565426"; "K"="1.035902087289256"; "L"="1.040889956875909"; "M"="1.048259905569686"; "N"="1.014965934352414" }
    @{ Row=15; "B"="1.02"; "C"="1.029586917491996"; "D"="1.032589274173725"; "E"="1.037581637873289"; "F"="1.045002393982244"; "I"="1.023915340507297"; "J"="1.036164389590152"; "K"="1.036153230001111"; "L"="1.041127137225676"; "M"="1.04852081925309"; "N"="1.015054155135528" }
    @{ Row=16; "B"="1.02"; "C"="1.031220381217136"; "D"="1.03416714562377"; "E"="1.03907731378174"; "F"="1.046635630581571"; "I"="1.023871434968932"; "J"="1.037570083988448"; "K"="1.037611203200613"; "L"="1.042503894736898"; "M"="1.050035662484549"; "N"="1.015565748614913" }
    @{ Row=17; "B"="1.02"; "C"="1.032241599264274"; "D"="1.035153676581915"; "E"="1.040012341989708"; "F"="1.047656888725776"; "I"="1.023843013862359"; "J"="1.038448548867654"; "K"="1.038522469008551"; "L"="1.043364258595995"; "M"="1.050982601642853"; "N"="1.015885011622206" }
    @{ Row=18; "B"="1.02"; "C"="1.032836052206622"; "D"="1.035727962332571"; "E"="1.040556606213861"; "F"="1.048251430112705"; "I"="1.023826118601775"; "J"="1.038959776764371"; "K"="1.039052832485775"; "L"="1.043864946101596"; "M"="1.051533773423066"; "N"="1.016070646239164" }
    @{ Row=19; "B"="1.02"; "C"="1.033038542831997"; "D"="1.035923587405972"; "E"="1.040741997864389"; "F"="1.048453961861148"; "I"="1.023820303822167"; "J"="1.039133896385259"; "K"="1.039233477464123"; "L"="1.044035474663893"; "M"="1.051721513654259"; "N"="1.016133844100684" }
    @{ Row=20; "B"="1.02"; "C"="1.032132157513226"; "D"="1.035047949654619"; "E"="1.039912138813007"; "F"="1.047547435927114"; "I"="1.023846096030143"; "J"="1.038354418929008"; "K"="1.038424819505771"; "L"="1.043272068890828"; "M"="1.050881124672924"; "N"="1.015850818504637" }
    @{ Row=21; "B"="1.02"; "C"="1.029171073216513"; "D"="1.032187603988182"; "E"="1.037200857333937"; "F"="1.044586664282698"; "I"="1.02392621804881"; "J"="1.035806420760361"; "K"="1.035781988357569"; "L"="1.04077653206984"; "M"="1.048135136650866"; "N"="1.014923736419071" }
    @{ Row=22; "B"="1.02"; "C"="1.027297961370661"; "D"="1.030378445588498"; "E"="1.03548561886827"; "F"="1.042714351047513"; "I"="1.023973741382756"; "J"="1.034193461159822"; "K"="1.034109422142072"; "L"="1.039196726786682"; "M"="1.046397705749013"; "N"="1.014335409212803" }
    @{ Row=23; "B"="1.02"; "C"="1.028292034507777"; "D"="1.031338556130767"; "E"="1.036395920750502"; "F"="1.043707942975485"; "I"="1.023948818288664"; "J"="1.035049578893487"; "K"="1.034997137271412"; "L"="1.040035252155886"; "M"="1.047319808003667"; "N"="1.014647815229793" }
    @{ Row=24; "B"="1.02"; "C"="1.032181613253358"; "D"="1.03509572662457"; "E"="1.039957419778833"; "F"="1.047596896458321"; "I"="1.023844704313072"; "J"="1.038396955795671"; "K"="1.038468946701523"; "L"="1.04331372899512"; "M"="1.050926981334269"; "N"="1.015866270712587" }
    @{ Row=25; "B"="1.02"; "C"="1.036648173342412"; "D"="1.039411178549705"; "E"="1.04404654526202"; "F"="1.052065275295968"; "I"="1.02371141742819"; "J"="1.042235876875733"; "K"="1.042452402005797"; "L"="1.047073381547944"; "M"="1.055067549748625"; "N"="1.017257303789772" }
)

foreach ($row in $data) {
    foreach ($col in @("B","C","D","E","F","I","J","K","L","M","N")) {
        $ws.Range($col + $row.Row).Value = [double]$row[$col]
    }
}